$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E (rows 2-51) to Text format so numeric-looking strings
# (e.g. "1.00", "113.70") are preserved verbatim instead of being coerced
# into numbers, matching the workbook's original inlineStr/text cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "52.103.87"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.843.60"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D5").Value = "362.62"
$ws.Range("E5").Value = "  +6.06%  "
$ws.Range("D6").Value = "113.70"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +4.07%  "
$ws.Range("D10").Value = "41.76"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "20.02"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "7.79"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "3.291.47"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "2.841.04"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "0.906"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "51.996.28"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +6.77%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").Value = "70.15"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "267.04"
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "27.11"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "10.42"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "53.55"
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").Value = "34.09"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("D34").Value = "0.0447"
$ws.Range("E34").Value = "  +20.95%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  +6.33%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0839"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").Value = "18.33"
$ws.Range("D41").Value = "24.31"
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D44").Value = "127.75"
$ws.Range("D45").Value = "2.26"
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("D46").Value = "2.121.75"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +9.93%  "
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("D51").Value = "9.02"
$ws.Range("E51").Value = "  +0.89%  "

# Restore default (unstyled) cell style now that values are committed as text,
# so the saved XML has no stray style index on these cells (matches original).
$dataRange.Style = "Normal"

